$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 26 (shifts the current rows 26-27 down to 27-28)
$ws.Rows(26).Insert()

# Populate the newly inserted row 26 with the updated weekly price record
$ws.Cells.Item(26, 1).Value = 3
$ws.Cells.Item(26, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(26, 3).Value = "Coquimbo"
$ws.Cells.Item(26, 4).Value = 44476
$ws.Cells.Item(26, 5).Value = 5
$ws.Cells.Item(26, 6).Value = 100112022
$ws.Cells.Item(26, 7).Value = "Arveja Verde"
$ws.Cells.Item(26, 8).Value = "Perfection"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 73
$ws.Cells.Item(26, 11).Value = 23000
$ws.Cells.Item(26, 12).Value = 24000
$ws.Cells.Item(26, 13).Value = 23521
$ws.Cells.Item(26, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(26, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(26, 16).Value = 941
$ws.Cells.Item(26, 17).Value = 25
$ws.Cells.Item(26, 18).Value = "Hortaliza"
